# Updates the cryptos list: refresh Price (D) and Volume(1h) (E) values,
# and fix two pairs of rows whose Coin/Link/Price/Volume data were swapped
# (OKB <-> ImmutableX at rows 45/46, dogwifhat <-> Mantle at rows 48/49).
#
# Price values are plain text in this sheet (e.g. "76.577.45", "1.00",
# "0.0000112"). Excel's COM layer auto-detects plain decimal-looking
# strings (a single dot, no thousands separators) as numbers, which would
# silently turn "1.00" into 1 or "5.10" into 5.1. To keep those as text we
# force NumberFormat="@" before writing on the cells where that risk
# applies. Multi-dot values (thousands grouping, e.g. "3.042.88") are never
# auto-parsed as numbers, so they don't need the text format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# row -> (D value, E value); $null means "leave unchanged".
# Entries marked Text=$true need the NumberFormat="@" guard.
$updates = @(
    @{ Row = 2;  D = "76.577.45";  DText = $false; E = "  +0.75%  " },
    @{ Row = 3;  D = "3.042.88";   DText = $false; E = "  +4.50%  " },
    @{ Row = 4;  D = $null;        DText = $false; E = "  -0.05%  " },
    @{ Row = 5;  D = "202.42";     DText = $true;  E = "  +0.06%  " },
    @{ Row = 6;  D = "630.82";     DText = $true;  E = "  +5.67%  " },
    @{ Row = 7;  D = "1.00";       DText = $true;  E = "  +0.02%  " },
    @{ Row = 8;  D = $null;        DText = $false; E = "  +0.73%  " },
    @{ Row = 9;  D = $null;        DText = $false; E = "  +6.50%  " },
    @{ Row = 10; D = "3.042.16";   DText = $false; E = "  +4.46%  " },
    @{ Row = 11; D = "0.438";      DText = $true;  E = "  +2.27%  " },
    @{ Row = 12; D = $null;        DText = $false; E = "  -0.43%  " },
    @{ Row = 13; D = "5.10";       DText = $true;  E = "  +4.87%  " },
    @{ Row = 14; D = "3.601.25";   DText = $false; E = "  +4.45%  " },
    @{ Row = 15; D = "29.53";      DText = $true;  E = "  +6.55%  " },
    @{ Row = 16; D = "76.598.49";  DText = $false; E = "  +0.92%  " },
    @{ Row = 17; D = $null;        DText = $false; E = "  +2.16%  " },
    @{ Row = 18; D = "3.046.23";   DText = $false; E = "  +4.67%  " },
    @{ Row = 19; D = "13.48";      DText = $true;  E = "  +4.27%  " },
    @{ Row = 20; D = "8.78";       DText = $true;  E = "  +0.50%  " },
    @{ Row = 21; D = "376.84";     DText = $true;  E = "  +1.02%  " },
    @{ Row = 22; D = $null;        DText = $false; E = "  +0.13%  " },
    @{ Row = 23; D = "4.37";       DText = $true;  E = "  +2.23%  " },
    @{ Row = 24; D = $null;        DText = $false; E = "  +4.10%  " },
    @{ Row = 26; D = "4.40";       DText = $true;  E = "  +4.59%  " },
    @{ Row = 27; D = "0.997";      DText = $true;  E = "  -0.22%  " },
    @{ Row = 28; D = $null;        DText = $false; E = "  +3.21%  " },
    @{ Row = 29; D = "0.0000112";  DText = $true;  E = "  +4.18%  " },
    @{ Row = 30; D = "1.00";       DText = $true;  E = "  +0.23%  " },
    @{ Row = 31; D = "8.34";       DText = $true;  E = "  +8.21%  " },
    @{ Row = 32; D = $null;        DText = $false; E = "  +1.60%  " },
    @{ Row = 33; D = "517.08";     DText = $true;  E = "  +2.87%  " },
    @{ Row = 34; D = "1.97";       DText = $true;  E = "  +8.67%  " },
    @{ Row = 35; D = $null;        DText = $false; E = "  -0.03%  " },
    @{ Row = 36; D = "20.89";      DText = $true;  E = $null },
    @{ Row = 37; D = $null;        DText = $false; E = "  -1.71%  " },
    @{ Row = 38; D = $null;        DText = $false; E = "  +10.61%  " },
    @{ Row = 39; D = $null;        DText = $false; E = "  +2.03%  " },
    @{ Row = 40; D = $null;        DText = $false; E = "  +3.76%  " },
    @{ Row = 41; D = "188.35";     DText = $true;  E = "  +4.16%  " },
    @{ Row = 42; D = "0.113";      DText = $true;  E = "  -0.34%  " },
    @{ Row = 43; D = $null;        DText = $false; E = "  +0.00%  " },
    @{ Row = 44; D = "5.23";       DText = $true;  E = "  +4.73%  " },
    @{ Row = 47; D = $null;        DText = $false; E = "  +1.59%  " },
    @{ Row = 50; D = $null;        DText = $false; E = "  +6.62%  " },
    @{ Row = 51; D = "3.89";       DText = $true;  E = "  +4.54%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        if ($u.DText) {
            Set-TextValue "D$r" $u.D
        } else {
            $ws.Range("D$r").Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = $u.E
    }
}

# Rows 45/46 swap: OKB <-> ImmutableX (with updated price/volume values)
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D45" "1.26"
$ws.Range("E45").Value = "  +6.04%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D46" "42.05"
$ws.Range("E46").Value = "  +4.74%  "

# Rows 48/49 swap: dogwifhat <-> Mantle (with updated price/volume values)
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D48" "0.731"
$ws.Range("E48").Value = "  +11.35%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "2.45"
$ws.Range("E49").Value = "  +4.36%  "
